# NV-29 Lâm Hoàng Phú 8-2024.xlsx
# - Insert a new worksheet "Đơn phụ phẫu 1" between "Đơn sale chính" and "Lương"
#   and populate it with the phụ phẫu 1 order detail + total row.
# - Update the "Lương" (salary) summary sheet with the computed totals that
#   flow from the new phụ phẫu 1 data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new sheet right after "Đơn sale chính" (i.e. before "Lương")
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$newWs = $wb.Worksheets.Add($null, $firstSheet)
$newWs.Name = "Đơn phụ phẫu 1"

# Header row
$newWs.Cells.Item(1,1).Value = "Tiền tố"
$newWs.Cells.Item(1,2).Value = "Mã dịch vụ"
$newWs.Cells.Item(1,3).Value = "Ngày thực hiện"
$newWs.Cells.Item(1,4).Value = "Cơ sở"
$newWs.Cells.Item(1,5).Value = "Khách hàng"
$newWs.Cells.Item(1,6).Value = "Nguồn khách"
$newWs.Cells.Item(1,7).Value = "Tên dịch vụ"
$newWs.Cells.Item(1,8).Value = "Phụ phẫu 1"
$newWs.Cells.Item(1,9).Value = "Công phụ phẫu 1"

# Data row (order detail)
$newWs.Cells.Item(2,1).Value = "HD-LUXURY"
$newWs.Cells.Item(2,2).Value = 614
$newWs.Cells.Item(2,3).NumberFormat = "@"
$newWs.Cells.Item(2,3).Value = "08-01-2024"
$newWs.Cells.Item(2,4).Value = "CẦN THƠ"
$newWs.Cells.Item(2,5).Value = "Trần Nguyễn Yến Linh"
$newWs.Cells.Item(2,6).Value = "Khách cũ"
$newWs.Cells.Item(2,7).Value = "Cắt mí"
$newWs.Cells.Item(2,8).Value = "Lâm Hoàng Phú"
$newWs.Cells.Item(2,9).Value = 50000

# Total row
$newWs.Cells.Item(3,1).Value = "Tổng"
$newWs.Cells.Item(3,2).Value = 1
$newWs.Cells.Item(3,9).Value = 50000

# ---------------------------------------------------------------------------
# 2. Update the "Lương" summary sheet with the new totals
# ---------------------------------------------------------------------------
$luong = $wb.Worksheets.Item("Lương")

$luong.Cells.Item(2,2).Value = 2                         # Tổng công tại CẦN THƠ
$luong.Cells.Item(3,2).Value = 285714.2857142857          # Lương cơ bản tại CẦN THƠ
$luong.Cells.Item(8,2).Value = 50000                      # Công phụ phẫu 1 tại CẦN THƠ
$luong.Cells.Item(31,2).Value = 335714.2857142857         # Tổng lương tại CẦN THƠ
$luong.Cells.Item(34,1).Value = "Tổng lương tại HỆ THỐNG"  # label update
$luong.Cells.Item(34,2).Value = 335714.2857142857          # Tổng lương
